$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.ClearFormats()
}

Set-TextValue 'D2' '35.326.35'
Set-TextValue 'E2' '  -0.06%  '
Set-TextValue 'D3' '1.912.08'
Set-TextValue 'E3' '  +0.21%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '0.726'
Set-TextValue 'E5' '  +8.85%  '
Set-TextValue 'D6' '255.54'
Set-TextValue 'E7' '  +0.05%  '
Set-TextValue 'D8' '42.24'
Set-TextValue 'E8' '  +1.42%  '
Set-TextValue 'E9' '  +5.40%  '
Set-TextValue 'D10' '53.27'
Set-TextValue 'E10' '  +0.14%  '
Set-TextValue 'E11' '  +7.07%  '
Set-TextValue 'E12' '  -0.43%  '
Set-TextValue 'D13' '13.13'
Set-TextValue 'E13' '  +6.70%  '
Set-TextValue 'D14' '2.190.07'
Set-TextValue 'E14' '  +0.17%  '
Set-TextValue 'E15' '  +5.58%  '
Set-TextValue 'E16' '  +4.21%  '
Set-TextValue 'D17' '1.914.25'
Set-TextValue 'E17' '  +0.40%  '
Set-TextValue 'D18' '35.306.19'
Set-TextValue 'E18' '  -0.21%  '
Set-TextValue 'D19' '75.10'
Set-TextValue 'E19' '  +4.10%  '
Set-TextValue 'D20' '0.0₃0850'
Set-TextValue 'E20' '  +3.69%  '
Set-TextValue 'E21' '  +1.90%  '
Set-TextValue 'D22' '13.18'
Set-TextValue 'E22' '  +5.61%  '
Set-TextValue 'E23' '  +7.18%  '
Set-TextValue 'E24' '  +0.08%  '
Set-TextValue 'E25' '  +7.42%  '
Set-TextValue 'E26' '  -0.26%  '
Set-TextValue 'D27' '167.10'
Set-TextValue 'E27' '  -2.36%  '
Set-TextValue 'D28' '8.82'
Set-TextValue 'E28' '  +4.37%  '
Set-TextValue 'D29' '18.85'
Set-TextValue 'E29' '  +2.65%  '
Set-TextValue 'E30' '  +4.29%  '
Set-TextValue 'D31' '4.128.98'
Set-TextValue 'E31' '  -0.65%  '
Set-TextValue 'E32' '  +26.22%  '
Set-TextValue 'E33' '  +5.38%  '
Set-TextValue 'E34' '  +15.79%  '
Set-TextValue 'E35' '  +4.89%  '
Set-TextValue 'E36' '  +4.65%  '
Set-TextValue 'E37' '  +0.03%  '
Set-TextValue 'D38' '0.927'
Set-TextValue 'E38' '  -1.71%  '
Set-TextValue 'D39' '2.04'
Set-TextValue 'E39' '  +0.28%  '
Set-TextValue 'D40' '100.46'
Set-TextValue 'E40' '  +11.73%  '
Set-TextValue 'E41' '  +6.27%  '
Set-TextValue 'D42' '17.06'
Set-TextValue 'E42' '  +4.90%  '
Set-TextValue 'D43' '1.13'
Set-TextValue 'E43' '  +2.52%  '
Set-TextValue 'D44' '0.0650'
Set-TextValue 'E44' '  -0.45%  '
Set-TextValue 'E45' '  +3.42%  '
Set-TextValue 'D46' '1.345.17'
Set-TextValue 'E46' '  +0.41%  '
Set-TextValue 'E47' '  +1.26%  '
Set-TextValue 'E48' '  +3.56%  '
Set-TextValue 'E49' '  -0.92%  '
Set-TextValue 'D50' '45.40'
Set-TextValue 'E50' '  -7.07%  '
Set-TextValue 'D51' '0.0760'
Set-TextValue 'E51' '  +7.33%  '
